$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "AST001"
$ws.Range("D3").Value = "AST001"
$ws.Range("D4").Value = "AST001"

$ws.Range("D3").Select()
